$d = $word.ActiveDocument

function Split-RunAtWord {
    param([string]$Anchor, [int]$PrefixLen, [int]$WordLen)
    $found = $d.Content.Duplicate
    $found.Find.Execute($Anchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $wordStart = $found.Start + $PrefixLen
    $wordEnd = $wordStart + $WordLen
    $sub = $d.Range($wordStart, $wordEnd)
    # Force Word to split the enclosing run around this sub-range by
    # toggling a character property on and back off again - this is the
    # same mechanism real Word uses to break one run into three runs
    # around a proofing-flagged word.
    $sub.Bold = 1
    $sub.Bold = 0
}

# 1) "Sanity check for inhibitory input ..." -> split out "check"
Split-RunAtWord "Sanity check for inhibitory" 7 5

# 2) "Same setup as Experiment C1 just for a longer run time." -> split out "run"
Split-RunAtWord "longer run time." 7 3

# 3) "Voltage increase by 8mV ..." -> split out "increase"
Split-RunAtWord "Voltage increase by 8mV" 8 8

# 4) "Final voltage increase by +- 27mV ..." -> split out "increase"
Split-RunAtWord "Final voltage increase by" 14 8

# 5) "Definitely need to increase the simulation time to get to a steady state.. possibly to 30 seconds." -> split out "state.."
Split-RunAtWord "a steady state.. possibly" 9 7

# 6) "In progress on Beast" (red text) -> "Complete" (default color)
$beast = $d.Content.Duplicate
$beast.Find.Execute("In progress on Beast", $true, $false, $false, $false, $false, $true, 1, $false, "Complete", 2)
